# Edit script for netherlands_tweede-divisie_2023-2024.xlsx
# Commit: "Atualizado por script em 02-12-2023 20:45"
#
# This update re-syncs the betting-odds sheet with the source scrape:
#  1) A handful of row-groups that share an identical kickoff timestamp (column E)
#     had their match-detail columns (F:home .. V:url_partida) shuffled between
#     rows by the scraper; this script restores the correct pairing.
#  2) Two newly played/scraped fixtures are appended as rows 121 and 122.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Permute F:V data among rows that share the same kickoff date/time (E column) ---
# Each group below represents rows whose match-detail columns (F..V) were
# shuffled in the source scrape; column A..E (index, country, league, season, datetime) stay put.

# Cycle: 2 -> 3 -> 2
$v2 = $ws.Range("F2:V2").Value()
$v3 = $ws.Range("F3:V3").Value()
$ws.Range("F2:V2").Value = $v3
$ws.Range("F3:V3").Value = $v2

# Cycle: 6 -> 7 -> 6
$v6 = $ws.Range("F6:V6").Value()
$v7 = $ws.Range("F7:V7").Value()
$ws.Range("F6:V6").Value = $v7
$ws.Range("F7:V7").Value = $v6

# Cycle: 18 -> 19 -> 20 -> 21 -> 18
$v18 = $ws.Range("F18:V18").Value()
$v19 = $ws.Range("F19:V19").Value()
$v20 = $ws.Range("F20:V20").Value()
$v21 = $ws.Range("F21:V21").Value()
$ws.Range("F18:V18").Value = $v19
$ws.Range("F19:V19").Value = $v20
$ws.Range("F20:V20").Value = $v21
$ws.Range("F21:V21").Value = $v18

# Cycle: 25 -> 27 -> 26 -> 25
$v25 = $ws.Range("F25:V25").Value()
$v27 = $ws.Range("F27:V27").Value()
$v26 = $ws.Range("F26:V26").Value()
$ws.Range("F25:V25").Value = $v27
$ws.Range("F27:V27").Value = $v26
$ws.Range("F26:V26").Value = $v25

# Cycle: 28 -> 29 -> 28
$v28 = $ws.Range("F28:V28").Value()
$v29 = $ws.Range("F29:V29").Value()
$ws.Range("F28:V28").Value = $v29
$ws.Range("F29:V29").Value = $v28

# Cycle: 33 -> 34 -> 35 -> 33
$v33 = $ws.Range("F33:V33").Value()
$v34 = $ws.Range("F34:V34").Value()
$v35 = $ws.Range("F35:V35").Value()
$ws.Range("F33:V33").Value = $v34
$ws.Range("F34:V34").Value = $v35
$ws.Range("F35:V35").Value = $v33

# Cycle: 59 -> 60 -> 59
$v59 = $ws.Range("F59:V59").Value()
$v60 = $ws.Range("F60:V60").Value()
$ws.Range("F59:V59").Value = $v60
$ws.Range("F60:V60").Value = $v59

# Cycle: 61 -> 64 -> 61
$v61 = $ws.Range("F61:V61").Value()
$v64 = $ws.Range("F64:V64").Value()
$ws.Range("F61:V61").Value = $v64
$ws.Range("F64:V64").Value = $v61

# Cycle: 62 -> 63 -> 62
$v62 = $ws.Range("F62:V62").Value()
$v63 = $ws.Range("F63:V63").Value()
$ws.Range("F62:V62").Value = $v63
$ws.Range("F63:V63").Value = $v62

# Cycle: 69 -> 70 -> 69
$v69 = $ws.Range("F69:V69").Value()
$v70 = $ws.Range("F70:V70").Value()
$ws.Range("F69:V69").Value = $v70
$ws.Range("F70:V70").Value = $v69

# Cycle: 72 -> 74 -> 72
$v72 = $ws.Range("F72:V72").Value()
$v74 = $ws.Range("F74:V74").Value()
$ws.Range("F72:V72").Value = $v74
$ws.Range("F74:V74").Value = $v72

# Cycle: 79 -> 80 -> 79
$v79 = $ws.Range("F79:V79").Value()
$v80 = $ws.Range("F80:V80").Value()
$ws.Range("F79:V79").Value = $v80
$ws.Range("F80:V80").Value = $v79

# Cycle: 98 -> 99 -> 98
$v98 = $ws.Range("F98:V98").Value()
$v99 = $ws.Range("F99:V99").Value()
$ws.Range("F98:V98").Value = $v99
$ws.Range("F99:V99").Value = $v98

# Cycle: 108 -> 109 -> 110 -> 108
$v108 = $ws.Range("F108:V108").Value()
$v109 = $ws.Range("F109:V109").Value()
$v110 = $ws.Range("F110:V110").Value()
$ws.Range("F108:V108").Value = $v109
$ws.Range("F109:V109").Value = $v110
$ws.Range("F110:V110").Value = $v108

# --- Append two new match rows (121, 122) discovered by this scrape run ---
# Row 121
$ws.Range("A120:V120").Copy($ws.Range("A121:V121"))
$ws.Range("A121").Value = 120
$ws.Range("B121").Value = 'netherlands'
$ws.Range("C121").Value = 'tweede-divisie'
$ws.Range("D121").Value = '2023-2024'
$ws.Range("E121").Value = 45262.625
$ws.Range("F121").Value = 'Quick Boys'
$ws.Range("G121").Value = 3
$ws.Range("H121").Value = 'Jong Sparta Rotterdam'
$ws.Range("I121").Value = 1
$ws.Range("J121").Value = 1.72
$ws.Range("K121").Value = '02/12/2023 03:12'
$ws.Range("L121").Value = 1.75
$ws.Range("M121").Value = '02/12/2023 14:58'
$ws.Range("N121").Value = 4.22
$ws.Range("O121").Value = '02/12/2023 03:12'
$ws.Range("P121").Value = 4.26
$ws.Range("Q121").Value = '02/12/2023 14:58'
$ws.Range("R121").Value = 3.65
$ws.Range("S121").Value = '02/12/2023 03:12'
$ws.Range("T121").Value = 3.57
$ws.Range("U121").Value = '02/12/2023 14:58'
$ws.Range("V121").Value = 'https://www.betexplorer.com/football/netherlands/tweede-divisie/quick-boys-jong-sparta-rotterdam/GSpIkBc4/'

# Row 122
$ws.Range("A121:V121").Copy($ws.Range("A122:V122"))
$ws.Range("A122").Value = 121
$ws.Range("B122").Value = 'netherlands'
$ws.Range("C122").Value = 'tweede-divisie'
$ws.Range("D122").Value = '2023-2024'
$ws.Range("E122").Value = 45262.64583333334
$ws.Range("F122").Value = 'Jong Almere City'
$ws.Range("G122").Value = 4
$ws.Range("H122").Value = 'Scheveningen'
$ws.Range("I122").Value = 1
$ws.Range("J122").Value = 2.42
$ws.Range("K122").Value = '02/12/2023 03:43'
$ws.Range("L122").Value = 2.25
$ws.Range("M122").Value = '02/12/2023 15:29'
$ws.Range("N122").Value = 3.49
$ws.Range("O122").Value = '02/12/2023 03:43'
$ws.Range("P122").Value = 3.62
$ws.Range("Q122").Value = '02/12/2023 15:29'
$ws.Range("R122").Value = 2.54
$ws.Range("S122").Value = '02/12/2023 03:43'
$ws.Range("T122").Value = 2.75
$ws.Range("U122").Value = '02/12/2023 15:29'
$ws.Range("V122").Value = 'https://www.betexplorer.com/football/netherlands/tweede-divisie/jong-almere-city-svv-scheveningen/UZ4NTTzo/'
